$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.319.84"
$ws.Range("E2").Value = "  +1.51%  "

$ws.Range("D3").Value = "3.151.87"
$ws.Range("E3").Value = "  +2.67%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'570.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.29%  "

$ws.Range("D6").Value = "'151.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.27%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "3.142.20"
$ws.Range("E8").Value = "  +2.53%  "

$ws.Range("E9").Value = "  +2.30%  "

$ws.Range("D10").Value = "'7.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +15.52%  "

$ws.Range("E11").Value = "  +1.93%  "

$ws.Range("D12").Value = "'0.472"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.81%  "

$ws.Range("D13").Value = "'36.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.36%  "

$ws.Range("E14").Value = "  +1.73%  "

$ws.Range("D15").Value = "3.656.57"
$ws.Range("E15").Value = "  +2.90%  "

$ws.Range("D16").Value = "65.331.40"
$ws.Range("E16").Value = "  +1.58%  "

$ws.Range("D17").Value = "'538.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +10.05%  "

$ws.Range("E18").Value = "  +2.17%  "

$ws.Range("D19").Value = "3.149.98"
$ws.Range("E19").Value = "  +2.75%  "

$ws.Range("D20").Value = "'6.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.72%  "

$ws.Range("D21").Value = "'14.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.49%  "

$ws.Range("D22").Value = "'0.716"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.87%  "

$ws.Range("D23").Value = "'7.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.81%  "

$ws.Range("D24").Value = "'13.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.44%  "

$ws.Range("D25").Value = "'79.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.91%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").Value = "'9.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.96%  "

$ws.Range("D28").Value = "'2.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.05%  "

$ws.Range("E29").Value = "  +5.19%  "

$ws.Range("D31").Value = "'2.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.67%  "

$ws.Range("D32").Value = "'26.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("D33").Value = "'1.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.12%  "

$ws.Range("D34").Value = "'556.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.90%  "

$ws.Range("D35").Value = "'5.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.22%  "

$ws.Range("D36").Value = "'6.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.95%  "

$ws.Range("D37").Value = "'0.0454"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.73%  "

$ws.Range("D38").Value = "'53.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.95%  "

$ws.Range("D39").Value = "'0.0835"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.33%  "

$ws.Range("D40").Value = "'2.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.79%  "

$ws.Range("D41").Value = "'0.124"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.12%  "

$ws.Range("D42").Value = "3.069.62"
$ws.Range("E42").Value = "  +6.31%  "

$ws.Range("D43").Value = "'8.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.04%  "

$ws.Range("E44").Value = "  +7.93%  "

$ws.Range("D45").Value = "'2.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.26%  "

$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("D47").Value = "'25.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.44%  "

$ws.Range("D48").Value = "0.0₃0536"
$ws.Range("E48").Value = "  -1.26%  "

$ws.Range("E49").Value = "  +2.80%  "

$ws.Range("D50").Value = "'120.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.32%  "

$ws.Range("D51").Value = "'2.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.90%  "
